$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-10 from 2023-10-13 (45212)
# to 2023-10-22 (45221), matching the automatic data refresh.
$ws.Range("C2:C10").Value = 45221
